# Insert a new "setup" worksheet as the 2nd sheet (right after
# "input_k_constants_log10"), make it the active sheet, and populate it
# with two rows of setup data: Calorimeter/DSC and Initial volume/15.

$wb = $excel.ActiveWorkbook

# Insert the new sheet before the current 2nd sheet ("input_concentrations")
# so it lands in tab position 2.
$setup = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$setup.Name = "setup"

# Populate the new sheet's data.
$setup.Range("A1").Value = "Calorimeter"
$setup.Range("B1").Value = "DSC"
$setup.Range("A2").Value = "Initial volume"
$setup.Range("B2").Value = 15

# Match the recorded selection/active cell on the new sheet.
$setup.Range("F6").Select()
